$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2-14
# from serial date 45174 (2023-09-05) to 45175 (2023-09-06)
for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value = 45175
    }
}
